# Rename existing sheet and add a new "Composições" sheet with cost
# composition records, mirroring the author's commit:
# "created records for composition and provided in a separate sheet"

$wb = $excel.ActiveWorkbook

# --- 1. Rename the first (existing) sheet ---------------------------------
$wsOrc = $wb.Worksheets.Item(1)
$wsOrc.Name = "Orçamento"

# --- 2. Insert a new sheet right after it, named "Composições" ------------
$wsComp = $wb.Worksheets.Add($null, $wsOrc)
$wsComp.Name = "Composições"

# --- 3. Fill in the cost-composition template ------------------------------
# Each cell is written individually (value + font) so that no incidental
# blank cells get materialised by range-wide formatting.

# Big title (C1)
$c = $wsComp.Range("C1")
$c.Value = "COMPOSIÇÕES DE CUSTO"
$c.Font.Size = 17
$c.Font.Bold = $true

# Sub-title row (C3, F3)
$c = $wsComp.Range("C3")
$c.Value = "COMPOSIÇÃO DE CUSTO UNITÁRIO"
$c.Font.Size = 13
$c.Font.Bold = $true

$c = $wsComp.Range("F3")
$c.Value = "comp_proj_date"
$c.Font.Size = 13
$c.Font.Bold = $true

# Field-label row (A5, C5, F5)
$c = $wsComp.Range("A5")
$c.Value = "comp_id"
$c.Font.Bold = $true

$c = $wsComp.Range("C5")
$c.Value = "comp_name"
$c.Font.Bold = $true

$c = $wsComp.Range("F5")
$c.Value = "comp_unit"
$c.Font.Bold = $true

# Table header row (A6:F6) - bold, teal
$headerCells = @("A6", "B6", "C6", "D6", "E6", "F6")
$headerValues = @("Código", "Descrição dos serviços", "Unid", "Coef", "R$ Unit", "R$ Total")
for ($i = 0; $i -lt $headerCells.Length; $i++) {
    $c = $wsComp.Range($headerCells[$i])
    $c.Value = $headerValues[$i]
    $c.Font.Bold = $true
    $c.Font.Color = 9400369
}

# Placeholder/merge-field row (A7:F7)
$fieldCells = @("A7", "B7", "C7", "D7", "E7", "F7")
$fieldValues = @("ID", "name", "unit", "bom_quant", "bom_avgprice", "bom_quant_avgprice")
for ($i = 0; $i -lt $fieldCells.Length; $i++) {
    $c = $wsComp.Range($fieldCells[$i])
    $c.Value = $fieldValues[$i]
    $c.Font.Size = 9
}

# Totals row (D8, F8)
$c = $wsComp.Range("D8")
$c.Value = "TOTAL MATERIAL/EQUIPAMENTO"
$c.Font.Size = 9

$c = $wsComp.Range("F8")
$c.Value = "all_total"
$c.Font.Size = 9

# --- 4. Column widths --------------------------------------------------
$wsComp.Columns.Item(2).ColumnWidth = 19.18
$wsComp.Columns.Item(3).ColumnWidth = 14.44
$wsComp.Columns.Item(4).ColumnWidth = 18.41
$wsComp.Columns.Item(5).ColumnWidth = 18.52
$wsComp.Columns.Item(6).ColumnWidth = 22.04

# --- 5. Row heights ------------------------------------------------------
$wsComp.Rows.Item(1).RowHeight = 20.9
$wsComp.Rows.Item(3).RowHeight = 12.75
$wsComp.Rows.Item(4).RowHeight = 12.75
$wsComp.Rows.Item(5).RowHeight = 12.75
$wsComp.Rows.Item(6).RowHeight = 12.75
$wsComp.Rows.Item(7).RowHeight = 12.75
$wsComp.Rows.Item(8).RowHeight = 12.75

# --- 6. Selection / active sheet -----------------------------------------
$wsOrc.Activate()
$wsOrc.Range("A6").Select()

$wsComp.Activate()
$wsComp.Range("F3").Select()
